# Update Name of Algo
# Apply the numeric corrections produced by re-running the RandomForest
# imputation algorithm on the terrestrial_mammals ACD/15/seed2 dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.45669999999999
$ws.Range("A3").Value = -21.39540000000003
$ws.Range("C5").Value = -14.3163
$ws.Range("D5").Value = -8.623099999999999
$ws.Range("D9").Value = -8.4557
$ws.Range("D11").Value = -8.389600000000005
$ws.Range("A14").Value = -20.40799999999998
$ws.Range("A21").Value = -21.20440000000001
$ws.Range("D21").Value = -7.900700000000006
$ws.Range("A23").Value = -21.73950000000004
$ws.Range("A25").Value = -22.47700000000003
